$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.082.58"

$ws.Range("D3").Value = "2.640.87"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.542"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.86%  "

$ws.Range("E9").Value = "  +1.01%  "

$ws.Range("E10").Value = "  -1.20%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.94"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.00%  "

$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").Value = "3.122.10"
$ws.Range("E15").Value = "  -0.33%  "

$ws.Range("D16").Value = "67.945.58"
$ws.Range("E16").Value = "  -1.20%  "

$ws.Range("D17").Value = "2.652.36"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "362.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.38"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.13%  "

$ws.Range("D27").Value = "2.774.80"
$ws.Range("E27").Value = "  -0.17%  "

$ws.Range("E28").Value = "  -2.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "553.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.96%  "

$ws.Range("E31").Value = "  -0.61%  "

$ws.Range("E32").Value = "  -1.87%  "

$ws.Range("E33").Value = "  -1.09%  "

$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.128"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.92%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").Value = "  -0.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.68%  "

$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("E40").Value = "  -3.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.93%  "

$ws.Range("D42").Value = "0.0₆0337"
$ws.Range("E42").Value = "  +5.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "158.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "21.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0783"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.14%  "

$ws.Range("E51").Value = "  -1.75%  "
